# ------------------------------------------------------------------
# CryCompanywiseStockReport_1.xlsx - stock count / value corrections
#
# The source system re-ran its stock valuation and produced corrected
# on-hand quantities (col F) and stock values (col G, = Rate(D) * Qty(F))
# for a number of items. Each company's "Sub Total:" row (col B) is the
# sum of the stock values of its items, and the final "Sub Total:" /
# "Grand Total:" rows (718/719) are the sum of every company subtotal,
# so those roll-up cells are corrected too.
#
# A handful of row-pairs (HIMALAYA WELLNESS COMPANY, HINDUSTAN UNILIVER
# LTD, Rasna Private Limited) also had their Code/Name/Qty/Value swapped
# between the two adjacent rows; those are applied below as direct value
# writes as well.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 3M INDIA LTD --
$ws.Range("F8").Value = 61
$ws.Range("G8").Value = 14199.58
# -- Sub Total / Grand Total roll-up --
$ws.Range("B9").Value = 18259.01
# -- ARAVIND LABORATORIES --
$ws.Range("F20").Value = 328
$ws.Range("G20").Value = 16803.44
$ws.Range("F24").Value = 25
$ws.Range("G24").Value = 2561.25
$ws.Range("F27").Value = 50
$ws.Range("G27").Value = 1793
$ws.Range("F29").Value = 70
$ws.Range("G29").Value = 3586.1
$ws.Range("F30").Value = 135
$ws.Range("G30").Value = 3804.3
$ws.Range("F32").Value = 30
$ws.Range("G32").Value = 1450.8
# -- Sub Total / Grand Total roll-up --
$ws.Range("B34").Value = 58255.73
# -- ASQUARE FOOD BEVERAGES PRIVATE LIMITED --
$ws.Range("F38").Value = 440
$ws.Range("G38").Value = 16020.4
$ws.Range("F45").Value = 83
$ws.Range("G45").Value = 1887.42
$ws.Range("F47").Value = 50
$ws.Range("G47").Value = 1264.5
$ws.Range("F48").Value = 235
$ws.Range("G48").Value = 13183.5
$ws.Range("F61").Value = 224
$ws.Range("G61").Value = 58403.52
# -- Sub Total / Grand Total roll-up --
$ws.Range("B66").Value = 202199.63
# -- Cholayil Pvt Ltd --
$ws.Range("F149").Value = 52
$ws.Range("G149").Value = 9427.08
# -- Sub Total / Grand Total roll-up --
$ws.Range("B155").Value = 36340.04
# -- DABUR INDIA LIMITED --
$ws.Range("F180").Value = 41
$ws.Range("G180").Value = 6906.04
$ws.Range("F182").Value = 17
$ws.Range("G182").Value = 1522.18
$ws.Range("F184").Value = 53
$ws.Range("G184").Value = 4346
# -- Sub Total / Grand Total roll-up --
$ws.Range("B193").Value = 63948.06
# -- GLAMIC HYGIENE PRODUCTS PVT LTD --
$ws.Range("F206").Value = 66
$ws.Range("G206").Value = 4276.8
# -- Sub Total / Grand Total roll-up --
$ws.Range("B208").Value = 4323.29
# -- Glaxosmithkline Asia Private Limited --
$ws.Range("F213").Value = 210
$ws.Range("G213").Value = 26602.8
$ws.Range("F217").Value = 34
$ws.Range("G217").Value = 2526.2
# -- Sub Total / Grand Total roll-up --
$ws.Range("B218").Value = 75910.08
# -- GODREJ CONSUMER PRODUCTS LIMITED --
$ws.Range("F222").Value = 791
$ws.Range("G222").Value = 14633.5
# -- Sub Total / Grand Total roll-up --
$ws.Range("B229").Value = 24521.33
# -- HIMALAYA WELLNESS COMPANY --
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F291").Value = 22
$ws.Range("G291").Value = 1884.96
$ws.Range("B292").Value = 66196
$ws.Range("C292").Value = 'HIM-Total Care Baby Pants Drapers-Xl-9S'
$ws.Range("F292").Value = 6
$ws.Range("G292").Value = 526.2
$ws.Range("B293").Value = 64985
$ws.Range("C293").Value = 'HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S'
$ws.Range("F293").Value = 12
$ws.Range("G293").Value = 1052.4
# -- HINDUSTAN UNILIVER LTD --
$ws.Range("B297").Value = 61610
$ws.Range("E297").Value = 122.71
$ws.Range("F297").Value = -58
$ws.Range("G297").Value = -5957.18
$ws.Range("B298").Value = 63565
$ws.Range("E298").Value = 109.19
$ws.Range("F298").Value = 60
$ws.Range("G298").Value = 6162.6
$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 27
$ws.Range("G306").Value = 3873.96
$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92
$ws.Range("F324").Value = 30
$ws.Range("G324").Value = 5139.9
$ws.Range("F325").Value = 26
$ws.Range("G325").Value = 3930.42
# -- Sub Total / Grand Total roll-up --
$ws.Range("B328").Value = -11041.45
# -- Kanav Enterprises Private Limited --
$ws.Range("F361").Value = 230
$ws.Range("G361").Value = 32335.7
# -- Sub Total / Grand Total roll-up --
$ws.Range("B363").Value = 74436.61
# -- KARNATAKA SOAPS & DETERGENTS LTD --
$ws.Range("F365").Value = 14
$ws.Range("G365").Value = 774.62
# -- Sub Total / Grand Total roll-up --
$ws.Range("B372").Value = 59194.65
# -- KUSHAL KARYASHALA PVT LTD --
$ws.Range("F387").Value = 419
$ws.Range("G387").Value = 40475.4
# -- Sub Total / Grand Total roll-up --
$ws.Range("B389").Value = 57357.38
# -- LIFE STYLE FOODS PVT LTD --
$ws.Range("F397").Value = 73
$ws.Range("G397").Value = 2640.41
$ws.Range("F408").Value = 205
$ws.Range("G408").Value = 3249.25
# -- Sub Total / Grand Total roll-up --
$ws.Range("B417").Value = 171009.76
# -- N.RANGA RAO & SONS PVT LTD --
$ws.Range("F432").Value = 93
$ws.Range("G432").Value = 4502.13
# -- Sub Total / Grand Total roll-up --
$ws.Range("B438").Value = 24735.15
# -- Rasna Private Limited --
$ws.Range("B506").Value = 60022
$ws.Range("E506").Value = 37.22
$ws.Range("F506").Value = -113
$ws.Range("G506").Value = -3709.79
$ws.Range("B507").Value = 64830
$ws.Range("E507").Value = 34.9
$ws.Range("F507").Value = 84
$ws.Range("G507").Value = 2757.72
# -- RECKITT BENCKISER INDIA PVT LTD --
$ws.Range("F511").Value = 228
$ws.Range("G511").Value = 22770.36
# -- Sub Total / Grand Total roll-up --
$ws.Range("B525").Value = 119637.15
# -- SARATHI INTERNATIONAL INC --
$ws.Range("F530").Value = 16
$ws.Range("G530").Value = 690.88
$ws.Range("F532").Value = 7
$ws.Range("G532").Value = 302.26
# -- Sub Total / Grand Total roll-up --
$ws.Range("B535").Value = 23592.02
# -- SOUTHERN HEALTH FOODS PVT LTD --
$ws.Range("F558").Value = 189
$ws.Range("G558").Value = 23029.65
# -- Sub Total / Grand Total roll-up --
$ws.Range("B561").Value = 27519.85
# -- SP APPARELS LTD --
$ws.Range("F563").Value = 18
$ws.Range("G563").Value = 3363.48
$ws.Range("F567").Value = 0
$ws.Range("G567").Value = 0
# -- Sub Total / Grand Total roll-up --
$ws.Range("B573").Value = 23237.65
# -- Tip Top Food Tech India --
$ws.Range("F614").Value = 79
$ws.Range("G614").Value = 11461.32
$ws.Range("F615").Value = 100
$ws.Range("G615").Value = 15467
# -- Sub Total / Grand Total roll-up --
$ws.Range("B628").Value = 209289.96
# -- V-GUARD INDUSTRIES LTD --
$ws.Range("F630").Value = 2
$ws.Range("G630").Value = 21707.02
# -- Sub Total / Grand Total roll-up --
$ws.Range("B635").Value = 111045.79
# -- Vishuddha Nutriol Private Limited --
$ws.Range("F662").Value = 39
$ws.Range("G662").Value = 3132.09
# -- Sub Total / Grand Total roll-up --
$ws.Range("B668").Value = 11991.81
# -- VVD AND SONS PRIVATE LIMITED --
$ws.Range("F674").Value = 754
$ws.Range("G674").Value = 122984.94
# -- Sub Total / Grand Total roll-up --
$ws.Range("B680").Value = 123997.49
# -- XO FOOTWEAR PVT LTD --
$ws.Range("F711").Value = 13
$ws.Range("G711").Value = 6940.05
# -- Sub Total / Grand Total roll-up --
$ws.Range("B713").Value = 65786.17
$ws.Range("B718").Value = 2642422.16
$ws.Range("B719").Value = 2642422.16
